$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("F1").Value = "Metadata - Single European Sky Portal"
$ws.Range("F2").Value = "pru-support@eurocontrol.int"
